$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I (shifts the existing "Cogm per kg" column from I to J)
$ws.Columns("I").Insert()

# New header for the inserted "Distribution channel code" column
$ws.Range("I1").Value = "Distribution channel code"

# New distribution channel code values for the two data rows
$ws.Range("I2").Value = "TR"
$ws.Range("I3").Value = "GO"

# Match the column width Excel computed for the new column
$ws.Columns("I").ColumnWidth = 21.72
